$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27
$ws.Cells.Item($row, 1).Value = 58
$ws.Cells.Item($row, 2).Value = "Update index.py"
$ws.Cells.Item($row, 3).Value = "riya-morankar"
$ws.Cells.Item($row, 4).Value = "riyam2309"
$ws.Cells.Item($row, 5).Value = "edit1 to main"
$dateCell = $ws.Cells.Item($row, 6)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-06-30"
$dateCell.Style = "Normal"
